$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values for the affected rows (trialTrain rows 1,8,11,19,21,27 -> sheet rows 2,9,12,20,22,28)
$ws.Range("D2").Value  = 5
$ws.Range("F2").Value  = 3
$ws.Range("H2").Value  = 46

$ws.Range("D9").Value  = 3
$ws.Range("F9").Value  = 3
$ws.Range("H9").Value  = 46

$ws.Range("D12").Value = 5
$ws.Range("F12").Value = 3
$ws.Range("H12").Value = 46

$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = 46

$ws.Range("D22").Value = 7
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 46

$ws.Range("D28").Value = 3
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 46

# Update the active selection to match the saved cursor position
$ws.Range("D28").Select()
